$data = New-Object 'object[,]' 13,11
$data[0,0] = "Codigo"
$data[0,1] = "Nivel1"
$data[0,2] = "Nivel2"
$data[0,3] = "Nivel3"
$data[0,4] = "Nivel4"
$data[0,5] = "Nivel5"
$data[0,6] = "Moneda"
$data[0,7] = "Es_Financiera"
$data[0,8] = "NivelCuenta"
$data[0,9] = "Fecha_Creacion"
$data[0,10] = "Fecha_Modificacion"
$data[1,0] = 1234567891
$data[1,1] = "Extra"
$data[1,2] = "Extra"
$data[1,3] = "Descuentos comerciales"
$data[1,4] = "Descuentos comerciales"
$data[1,5] = "Descuentos comerciales"
$data[1,6] = "##"
$data[1,7] = "N"
$data[1,8] = 5
$data[1,9] = 45292
$data[1,10] = 45292
$data[2,0] = 1234567892
$data[2,1] = "Extra"
$data[2,2] = "Extra"
$data[2,3] = "Venta a precio público (sin iva)"
$data[2,4] = "Venta a precio público (sin iva)"
$data[2,5] = "Venta a precio público (sin iva)"
$data[2,6] = "##"
$data[2,7] = "N"
$data[2,8] = 5
$data[2,9] = 45292
$data[2,10] = 45292
$data[3,0] = 1234567893
$data[3,1] = "Extra"
$data[3,2] = "Extra"
$data[3,3] = "Unidades"
$data[3,4] = "Unidades"
$data[3,5] = "Unidades"
$data[3,6] = "##"
$data[3,7] = "N"
$data[3,8] = 5
$data[3,9] = 45292
$data[3,10] = 45292
$data[4,0] = 4106020001
$data[4,1] = "Resultado"
$data[4,2] = "MARGEN"
$data[4,3] = "Mermas y Diferencia Inventarios"
$data[4,4] = "Mermas/Dif Inventario"
$data[4,5] = "Obsolescencia"
$data[4,6] = "##"
$data[4,7] = "N"
$data[4,8] = 5
$data[4,9] = 45292
$data[4,10] = 45292
$data[5,0] = 99997
$data[5,1] = "Patrimonio Neto"
$data[5,2] = "RESULTADOS"
$data[5,3] = "Resultados no asignados"
$data[5,4] = "Resultados no aignados"
$data[5,5] = "Resultados no asignados"
$data[5,6] = "##"
$data[5,7] = "N"
$data[5,8] = 5
$data[5,9] = 45292
$data[5,10] = 45292
$data[6,0] = 99999
$data[6,1] = "Patrimonio Neto"
$data[6,2] = "RESULTADOS"
$data[6,3] = "Resultados no asignados"
$data[6,4] = "Resultados no aignados"
$data[6,5] = "Resultados no asignados"
$data[6,6] = "##"
$data[6,7] = "N"
$data[6,8] = 5
$data[6,9] = 45292
$data[6,10] = 45292
$data[7,0] = 99990
$data[7,1] = "Otros"
$data[7,2] = "Otros"
$data[7,3] = "Otros"
$data[7,4] = "Otros"
$data[7,5] = "Otros"
$data[7,6] = "##"
$data[7,7] = "N"
$data[7,8] = 5
$data[7,9] = 45292
$data[7,10] = 45292
$data[8,0] = 99993
$data[8,1] = "Otros"
$data[8,2] = "Otros"
$data[8,3] = "Otros"
$data[8,4] = "Otros"
$data[8,5] = "Otros"
$data[8,6] = "##"
$data[8,7] = "N"
$data[8,8] = 5
$data[8,9] = 45292
$data[8,10] = 45292
$data[9,0] = 99994
$data[9,1] = "Otros"
$data[9,2] = "Otros"
$data[9,3] = "Otros"
$data[9,4] = "Otros"
$data[9,5] = "Otros"
$data[9,6] = "##"
$data[9,7] = "N"
$data[9,8] = 5
$data[9,9] = 45292
$data[9,10] = 45292
$data[10,0] = 99995
$data[10,1] = "Otros"
$data[10,2] = "Otros"
$data[10,3] = "Otros"
$data[10,4] = "Otros"
$data[10,5] = "Otros"
$data[10,6] = "##"
$data[10,7] = "N"
$data[10,8] = 5
$data[10,9] = 45292
$data[10,10] = 45292
$data[11,0] = 99996
$data[11,1] = "Otros"
$data[11,2] = "Otros"
$data[11,3] = "Otros"
$data[11,4] = "Otros"
$data[11,5] = "Otros"
$data[11,6] = "##"
$data[11,7] = "N"
$data[11,8] = 5
$data[11,9] = 45292
$data[11,10] = 45292
$data[12,0] = 99998
$data[12,1] = "Otros"
$data[12,2] = "Otros"
$data[12,3] = "Otros"
$data[12,4] = "Otros"
$data[12,5] = "Otros"
$data[12,6] = "##"
$data[12,7] = "N"
$data[12,8] = 5
$data[12,9] = 45292
$data[12,10] = 45292
$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook
$old = $wb.Worksheets.Item(1)
$new = $wb.Worksheets.Add()
$new.Range("A1:K13").Value = $data
$old.Delete()
$new.Name = "Hoja1"
$new.Range("J1").Value = "FechaCreacion"
$new.Range("K1").Value = "FechaModificacion"
Write-Output "done"
